$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Add new column K (2022 data), copying the formatting from the corresponding J cell in each row

$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial($xlPasteFormats)
$ws.Range("K4").Value = 2022

$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial($xlPasteFormats)
$ws.Range("K5").Value = 1.6

$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial($xlPasteFormats)
$ws.Range("K6").Value = 0.4

$ws.Range("J7").Copy()
$ws.Range("K7").PasteSpecial($xlPasteFormats)
$ws.Range("K7").Value = 0.9

$ws.Range("J8").Copy()
$ws.Range("K8").PasteSpecial($xlPasteFormats)
$ws.Range("K8").Value = 0.6

$ws.Range("J9").Copy()
$ws.Range("K9").PasteSpecial($xlPasteFormats)
$ws.Range("K9").Value = 2.1

$ws.Range("J10").Copy()
$ws.Range("K10").PasteSpecial($xlPasteFormats)
$ws.Range("K10").Value = 0.6

$ws.Range("J11").Copy()
$ws.Range("K11").PasteSpecial($xlPasteFormats)
$ws.Range("K11").Value = 0.9

$ws.Range("J12").Copy()
$ws.Range("K12").PasteSpecial($xlPasteFormats)
$ws.Range("K12").Value = 2.2999999999999998

$ws.Range("J13").Copy()
$ws.Range("K13").PasteSpecial($xlPasteFormats)
$ws.Range("K13").Value = 4.3

$ws.Range("J14").Copy()
$ws.Range("K14").PasteSpecial($xlPasteFormats)
$ws.Range("K14").Value = 0.3

$excel.CutCopyMode = $false

# Update selection to match the new active cell seen in the diff
$ws.Range("L7").Select()
